$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 196 (shifts existing rows 196:302 down to 197:303)
$ws.Rows("196:196").Insert()

# Populate the newly inserted row 196 with its data
$ws.Cells.Item(196, 1).Value = 4
$ws.Cells.Item(196, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(196, 3).Value = "Los Lagos"
$ws.Cells.Item(196, 4).Value = 44529
$ws.Cells.Item(196, 5).Value = 10
$ws.Cells.Item(196, 6).Value = "Fruta"
$ws.Cells.Item(196, 7).Value = 100102
$ws.Cells.Item(196, 8).Value = "Cítricos"
$ws.Cells.Item(196, 9).Value = 100102005
$ws.Cells.Item(196, 10).Value = "Naranja"
$ws.Cells.Item(196, 11).Value = "Navel Late"
$ws.Cells.Item(196, 12).Value = "Primera"
$ws.Cells.Item(196, 13).Value = 200
$ws.Cells.Item(196, 14).Value = 14000
$ws.Cells.Item(196, 15).Value = 14500
$ws.Cells.Item(196, 16).Value = 14250
$ws.Cells.Item(196, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(196, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(196, 19).Value = 950
$ws.Cells.Item(196, 20).Value = 15
